$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "57 x 26" + [char]11 + "  2    6" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "7|    |"
$t.Cell(1, 2).Range.Text = "33 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "3|    |"
$t.Cell(1, 3).Range.Text = "49 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "9|    |"
$t.Cell(2, 1).Range.Text = "66 x 55" + [char]11 + "  5    5" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(2, 2).Range.Text = "14 x 24" + [char]11 + "  2    4" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "4|    |"
$t.Cell(2, 3).Range.Text = "84 x 42" + [char]11 + "  4    2" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "4|    |"
$t.Cell(3, 1).Range.Text = "15 x 48" + [char]11 + "  4    8" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "5|    |"
$t.Cell(3, 2).Range.Text = "18 x 75" + [char]11 + "  7    5" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "8|    |"
$t.Cell(3, 3).Range.Text = "61 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "1|    |"
$t.Cell(4, 1).Range.Text = "42 x 41" + [char]11 + "  4    1" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "2|    |"
$t.Cell(4, 2).Range.Text = "73 x 48" + [char]11 + "  4    8" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "3|    |"
$t.Cell(4, 3).Range.Text = "47 x 24" + [char]11 + "  2    4" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "7|    |"
$t.Cell(5, 1).Range.Text = "96 x 45" + [char]11 + "  4    5" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "6|    |"
$t.Cell(5, 2).Range.Text = "67 x 66" + [char]11 + "  6    6" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "7|    |"
$t.Cell(5, 3).Range.Text = "26 x 17" + [char]11 + "  1    7" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "6|    |"
